$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")

# Mark attendance ("P") for the students/dates that were previously left
# blank in the Chamada (roll call) sheet.
$ws.Range("AQ5").Value = "P"
$ws.Range("AP6").Value = "P"
$ws.Range("AQ6").Value = "P"
$ws.Range("AP12").Value = "P"
$ws.Range("AQ12").Value = "P"
